$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rolling "Ultimo" timestamp update: every row block shifts up to the
# previous block's value, and the topmost block (rows 2-15) gets the
# fresh timestamp from this run.
for ($r = 2; $r -le 15; $r++) {
    $ws.Cells.Item($r, 4).Value = 44260.54314747927
}
for ($r = 16; $r -le 29; $r++) {
    $ws.Cells.Item($r, 4).Value = 44260.52179451389
}
for ($r = 30; $r -le 43; $r++) {
    $ws.Cells.Item($r, 4).Value = 44260.50046984954
}
